$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows at row 291 (pushing the existing rows 291-398 down to 293-400)
$ws.Rows("291:292").Insert()

# New row 291: weekly price record for Especial quality, dated 2022-11-11 (serial 44876)
$ws.Range("A291").Value = 11
$ws.Range("B291").Value = "Vega Monumental Concepción"
$ws.Range("C291").Value = "Bíobío"
$ws.Range("D291").Value = 44876
$ws.Range("E291").Value = 8
$ws.Range("F291").Value = "Fruta"
$ws.Range("G291").Value = 100101
$ws.Range("H291").Value = "Berries"
$ws.Range("I291").Value = 100112025
$ws.Range("J291").Value = "Frutilla"
$ws.Range("K291").Value = "Sin especificar"
$ws.Range("L291").Value = "Especial"
$ws.Range("M291").Value = 150
$ws.Range("N291").Value = 7000
$ws.Range("O291").Value = 7000
$ws.Range("P291").Value = 7000
$ws.Range("Q291").Value = "$/bandeja 7 kilos"
$ws.Range("R291").Value = "Provincia de Melipilla"
$ws.Range("S291").Value = 1000
$ws.Range("T291").Value = 7

# New row 292: weekly price record for Primera quality, dated 2022-11-11 (serial 44876)
$ws.Range("A292").Value = 11
$ws.Range("B292").Value = "Vega Monumental Concepción"
$ws.Range("C292").Value = "Bíobío"
$ws.Range("D292").Value = 44876
$ws.Range("E292").Value = 8
$ws.Range("F292").Value = "Fruta"
$ws.Range("G292").Value = 100101
$ws.Range("H292").Value = "Berries"
$ws.Range("I292").Value = 100112025
$ws.Range("J292").Value = "Frutilla"
$ws.Range("K292").Value = "Sin especificar"
$ws.Range("L292").Value = "Primera"
$ws.Range("M292").Value = 320
$ws.Range("N292").Value = 6000
$ws.Range("O292").Value = 6500
$ws.Range("P292").Value = 6188
$ws.Range("Q292").Value = "$/bandeja 7 kilos"
$ws.Range("R292").Value = "Provincia de Melipilla"
$ws.Range("S292").Value = 884
$ws.Range("T292").Value = 7
